$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (Volume/Number and date range) ----
$ws.Range("A8").Value = "Volume 31   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/28/2024  Through  11/3/2024"

# ---- Cells that change numeric<->text "N/A" style; copy a same-styled donor cell first ----
# Donor cells (untouched row 14): C14 -> text "0"; E14 -> text "***.* "; I14 -> count style; K14 -> pct style
$ws.Range("I14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("C23"))
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("C14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))
$ws.Range("I14").Copy($ws.Range("C33"))
$ws.Range("C33").Value = 1
$ws.Range("I14").Copy($ws.Range("F33"))
$ws.Range("F33").Value = 1
$ws.Range("I14").Copy($ws.Range("I33"))
$ws.Range("I33").Value = 1

# ---- Remaining cells: value-only updates (style unchanged) ----
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 9
$ws.Range("K15").Value = -55.555555555555
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 81
$ws.Range("J16").Value = 110
$ws.Range("K16").Value = -26.363636363636
$ws.Range("L16").Value = -20.588235294117
$ws.Range("M16").Value = -31.932773109243
$ws.Range("N16").Value = -88.964577656675
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 133
$ws.Range("J17").Value = 159
$ws.Range("K17").Value = -16.352201257861
$ws.Range("L17").Value = 7.258064516129
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = -40.090090090090
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 18.181818181818
$ws.Range("I18").Value = 116
$ws.Range("K18").Value = -29.268292682926
$ws.Range("L18").Value = -3.333333333333
$ws.Range("M18").Value = -23.684210526315
$ws.Range("N18").Value = -82.634730538922
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -18.75
$ws.Range("I19").Value = 373
$ws.Range("J19").Value = 490
$ws.Range("K19").Value = -23.877551020408
$ws.Range("L19").Value = -31.433823529411
$ws.Range("M19").Value = 8.115942028985
$ws.Range("N19").Value = -12.235294117647
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -81.818181818181
$ws.Range("J20").Value = 108
$ws.Range("K20").Value = -18.518518518518
$ws.Range("L20").Value = -22.123893805309
$ws.Range("M20").Value = 12.820512820512
$ws.Range("N20").Value = -89.873417721519
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -40
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = -21.978021978022
$ws.Range("I21").Value = 797
$ws.Range("J21").Value = 1042
$ws.Range("K21").Value = -23.512476007677
$ws.Range("L21").Value = -21.709233791748
$ws.Range("M21").Value = 8.141112618724
$ws.Range("N21").Value = -72.946367956551
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 26
$ws.Range("K22").Value = -19.230769230769
$ws.Range("L22").Value = -16
$ws.Range("L23").Value = 11.111111111111
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 57.142857142857
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = -17.977528089887
$ws.Range("I24").Value = 975
$ws.Range("J24").Value = 1221
$ws.Range("K24").Value = -20.147420147420
$ws.Range("L24").Value = -15.291051259774
$ws.Range("M24").Value = 24.203821656051
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 27.272727272727
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 67
$ws.Range("H25").Value = -31.343283582089
$ws.Range("I25").Value = 620
$ws.Range("J25").Value = 889
$ws.Range("K25").Value = -30.258717660292
$ws.Range("L25").Value = -27.230046948356
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = -21.052631578947
$ws.Range("I26").Value = 211
$ws.Range("J26").Value = 255
$ws.Range("K26").Value = -17.254901960784
$ws.Range("L26").Value = -18.846153846153
$ws.Range("M26").Value = 33.544303797468
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = -35.294117647058
$ws.Range("L27").Value = -31.25
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -33.333333333333
$ws.Range("N29").Value = -58.333333333333
$ws.Range("N30").Value = -58.333333333333
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
